# Apply the "Generate Report for Handback" change:
#  - Overview status cells: "Ready for handoff" -> "Handed back: in sync with en-US"
#  - zh-cn / de-de detail sheets: fill in "Latest Target File" (I) and
#    "Latest Handback File" (J) columns with hyperlinked .md / .xlf file
#    names, and stamp the "Latest Handback DateTime" (K) column.
#  - Widen columns that now hold the longer values.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$mdTarget964 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7f27b66b58c4f02ee7d201604ebebaed7ea59dc7/e2e/964f91e3-1ac1-4c3e-a233-ead2b28318ad.md"
$mdTargetB2d = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7f27b66b58c4f02ee7d201604ebebaed7ea59dc7/e2e/b2d59ada-7496-419c-a019-e9a56bb05d9f.md"
$mdName964 = "964f91e3-1ac1-4c3e-a233-ead2b28318ad.md"
$mdNameB2d = "b2d59ada-7496-419c-a019-e9a56bb05d9f.md"

# ---------------------------------------------------------------------
# Overview sheet: status text for both locale columns
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column also reflects the new text
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

# Row 2 (964f91e3 file)
$wsZh.Range("J2").Value = "964f91e3-1ac1-4c3e-a233-ead2b28318ad.44d0566a637d56e9f48ed31cfb94e60c0e650c64.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-18 16:27:48"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdTarget964, $null, $null, $mdName964)

# Row 3 (b2d59ada file)
$wsZh.Range("J3").Value = "b2d59ada-7496-419c-a019-e9a56bb05d9f.85225337f9082890dece64963425dc4b0bd8af34.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-18 16:27:48"
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdTargetB2d, $null, $null, $mdNameB2d)

$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Status column also reflects the new text
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# Row 2 (964f91e3 file)
$wsDe.Range("J2").Value = "964f91e3-1ac1-4c3e-a233-ead2b28318ad.44d0566a637d56e9f48ed31cfb94e60c0e650c64.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-18 16:27:54"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdTarget964, $null, $null, $mdName964)

# Row 3 (b2d59ada file)
$wsDe.Range("J3").Value = "b2d59ada-7496-419c-a019-e9a56bb05d9f.85225337f9082890dece64963425dc4b0bd8af34.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-18 16:27:54"
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdTargetB2d, $null, $null, $mdNameB2d)

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664

Write-Host "Report generated for handback."
